# Corrección de errores 1.0
# Removes the redundant "poder" from several "El sistema debe poder ..."
# requirement sentences, relocates the stray "_GoBack" bookmark to the
# end of the list (where Word last left the cursor), and adds a simple
# page-number footer to the section.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "El sistema debe poder mostrar al cliente un número correlativo de reclamo" `
             "El sistema debe mostrar al cliente un número correlativo de reclamo"

Replace-Text "El sistema debe poder filtrar los reclamos del cliente según la fecha." `
             "El sistema debe filtrar los reclamos del cliente según la fecha."

Replace-Text "El sistema debe poder filtrar los reclamos del cliente según el cuidador. " `
             "El sistema debe filtrar los reclamos del cliente según el cuidador. "

Replace-Text "El sistema debe poder filtrar los reclamos del cliente según el estado." `
             "El sistema debe filtrar los reclamos del cliente según el estado."

Replace-Text "El sistema debe poder filtrar un reclamo por su número correlativo." `
             "El sistema debe filtrar un reclamo por su número correlativo."

Replace-Text "El sistema debe poder mostrar un listado de todos los reclamos." `
             "El sistema debe  mostrar un listado de todos los reclamos."

Replace-Text "El sistema debe poder permitir la modificación del estado del reclamo." `
             "El sistema debe permitir la modificación del estado del reclamo."

Replace-Text "El sistema debe poder almacenar los informes de reclamo." `
             "El sistema debe almacenar los informes de reclamo."

Replace-Text "El sistema debe poder permitir la modificación de un informe de reclamo." `
             "El sistema debe permitir la modificación de un informe de reclamo."

# Insert a simple right-aligned page-number footer (Insert > Page Number >
# Bottom of Page), which is what leaves the "_GoBack" bookmark at the very
# end of the document content (the last place edited).
$sec = $d.Sections(1)
$footer = $sec.Footers(1)
$footer.PageNumbers.Add(2, $true) | Out-Null

Write-Output "done"
